# Insert a new data row at row 1119 (pushing all subsequent rows down by one),
# and populate it with a new weekly price record for Tomate - Primera.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1119:1119").Insert()

$ws.Range("A1119").Value = 4
$ws.Range("B1119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C1119").Value = "Los Lagos"
$ws.Range("D1119").Value = 45147
$ws.Range("E1119").Value = 10
$ws.Range("F1119").Value = 100112020
$ws.Range("G1119").Value = "Tomate"
$ws.Range("H1119").Value = "Larga vida"
$ws.Range("I1119").Value = "Primera"
$ws.Range("J1119").Value = 120
$ws.Range("K1119").Value = 27000
$ws.Range("L1119").Value = 27000
$ws.Range("M1119").Value = 27000
$ws.Range("N1119").Value = "$/bandeja 18 kilos"
$ws.Range("O1119").Value = "Región de Arica y Parinacota"
$ws.Range("P1119").Value = 1500
$ws.Range("Q1119").Value = 18
$ws.Range("R1119").Value = "Hortaliza"
